# New crime data collected - weekly update for cs-en-us-104pct.xlsx
# Applies:
#   1. Header text updates (volume number, week-covering dates) inside
#      rich-text shared strings.
#   2. Refreshed weekly/28-day/YTD/2yr crime statistics for rows 14-30,
#      including a handful of cells that flip between numeric "0 complaints"
#      values and the sheet's textual placeholders ("0" / "***.*" shown as
#      text where the underlying metric is undefined, e.g. 0/0 % change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

function Set-TextCell {
    param(
        [string]$Target,
        [string]$Text,
        [string]$FormatDonor
    )
    $dst = $ws.Range($Target)
    $dst.NumberFormat = "@"
    $dst.Value = $Text
    $src = $ws.Range($FormatDonor)
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

function Set-NumberCell {
    param(
        [string]$Target,
        $Number,
        [string]$FormatDonor
    )
    $dst = $ws.Range($Target)
    $dst.Value = $Number
    $src = $ws.Range($FormatDonor)
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# 1) Cells that change TYPE (number <-> text placeholder) this week.
#    Donors: row 23 is an all-text placeholder row (style for "text" cells);
#    D14/E14 stay text this week too (style for the other text donor);
#    row 16's numeric cells keep their normal numeric styles.
# ---------------------------------------------------------------------

# number -> text "0"
Set-TextCell -Target "C14" -Text "0" -FormatDonor "C23"
Set-TextCell -Target "G14" -Text "0" -FormatDonor "G23"
Set-TextCell -Target "F15" -Text "0" -FormatDonor "F23"
Set-TextCell -Target "C22" -Text "0" -FormatDonor "C23"
Set-TextCell -Target "D22" -Text "0" -FormatDonor "D23"
Set-TextCell -Target "C27" -Text "0" -FormatDonor "C23"
Set-TextCell -Target "C28" -Text "0" -FormatDonor "C23"
Set-TextCell -Target "C29" -Text "0" -FormatDonor "C23"
Set-TextCell -Target "D30" -Text "0" -FormatDonor "D23"

# number -> text "***.*"
Set-TextCell -Target "H14" -Text "***.*" -FormatDonor "H23"
Set-TextCell -Target "E22" -Text "***.*" -FormatDonor "E23"
Set-TextCell -Target "E30" -Text "***.*" -FormatDonor "E23"

# text -> number
Set-NumberCell -Target "D26" -Number 1 -FormatDonor "D16"
Set-NumberCell -Target "E26" -Number -100 -FormatDonor "E16"
Set-NumberCell -Target "D27" -Number 3 -FormatDonor "D16"
Set-NumberCell -Target "E27" -Number -100 -FormatDonor "E16"

# ---------------------------------------------------------------------
# 2) Pure value refreshes (style/type unchanged).
# ---------------------------------------------------------------------

$ws.Range("M14").Value = 0

$ws.Range("H15").Value = -100
$ws.Range("N15").Value = -5

$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 196
$ws.Range("J16").Value = 169
$ws.Range("K16").Value = 15.976331360946
$ws.Range("L16").Value = 78.181818181818
$ws.Range("M16").Value = -11.312217194570
$ws.Range("N16").Value = -76.039119804401

$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -22.222222222222
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -18.518518518518
$ws.Range("I17").Value = 233
$ws.Range("J17").Value = 252
$ws.Range("K17").Value = -7.539682539682
$ws.Range("L17").Value = 11.483253588516
$ws.Range("M17").Value = 28.021978021978
$ws.Range("N17").Value = -8.267716535433

$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -4.545454545454
$ws.Range("I18").Value = 168
$ws.Range("J18").Value = 231
$ws.Range("K18").Value = -27.272727272727
$ws.Range("L18").Value = -7.182320441988
$ws.Range("M18").Value = -55.319148936170
$ws.Range("N18").Value = -89.928057553956

$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 45.454545454545
$ws.Range("F19").Value = 69
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = 32.692307692307
$ws.Range("I19").Value = 577
$ws.Range("J19").Value = 543
$ws.Range("K19").Value = 6.261510128913
$ws.Range("L19").Value = 30.839002267573
$ws.Range("M19").Value = 59.833795013850
$ws.Range("N19").Value = 8.662900188323

$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -14.285714285714
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 39
$ws.Range("H20").Value = -41.025641025641
$ws.Range("I20").Value = 317
$ws.Range("J20").Value = 255
$ws.Range("K20").Value = 24.313725490196
$ws.Range("L20").Value = 83.236994219653
$ws.Range("M20").Value = 2.258064516129
$ws.Range("N20").Value = -89.180887372013

$ws.Range("C21").Value = 39
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 160
$ws.Range("G21").Value = 157
$ws.Range("H21").Value = 1.910828025477
$ws.Range("I21").Value = 1514
$ws.Range("J21").Value = 1471
$ws.Range("K21").Value = 2.923181509177
$ws.Range("L21").Value = 34.219858156028
$ws.Range("M21").Value = 2.923181509177
$ws.Range("N21").Value = -75.729400448861

$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666

$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -24.137931034482
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 122
$ws.Range("H24").Value = -13.934426229508
$ws.Range("I24").Value = 1103
$ws.Range("J24").Value = 1222
$ws.Range("K24").Value = -9.738134206219
$ws.Range("L24").Value = -0.541027953110
$ws.Range("M24").Value = 24.915062287655

$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 87.5
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = -5.405405405405
$ws.Range("I25").Value = 393
$ws.Range("J25").Value = 445
$ws.Range("K25").Value = -11.685393258427
$ws.Range("L25").Value = -6.650831353919
$ws.Range("M25").Value = -34.608985024958

$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 28
$ws.Range("K26").Value = 7.142857142857

$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 63
$ws.Range("K27").Value = -15.873015873015
$ws.Range("L27").Value = 17.777777777777

$ws.Range("M28").Value = 50

$ws.Range("M29").Value = 66.666666666666

# ---------------------------------------------------------------------
# 3) Header text: bump the volume/number and the reporting week dates.
#    These are rich-text shared strings; re-assigning .Characters(...)
#    .Text edits just the affected run's text in place.
# ---------------------------------------------------------------------

$volCell = $ws.Range("A8")
$volText = $volCell.Text
$numStart = $volText.Length - 1
$volCell.Characters($numStart, 2).Text = "43"

$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 10).Text = "10/23/2023"
$weekCell.Characters(48, 10).Text = "10/29/2023"
